$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 46015
$ws.Range("D8").Value = 152.34
$ws.Range("E8").Value = 152.05000000000001
$ws.Range("F8").Value = 162.05000000000001
$ws.Range("G8").Value = 152.16999999999999

# Row 9
$ws.Range("A9").Value = 46015
$ws.Range("D9").Value = 152.34
$ws.Range("E9").Value = 152.05000000000001
$ws.Range("F9").Value = 162.05000000000001
$ws.Range("G9").Value = 152.16999999999999

# Row 10
$ws.Range("A10").Value = 46015
$ws.Range("D10").Value = 155.07
$ws.Range("E10").Value = 154.18
$ws.Range("F10").Value = 164.18
$ws.Range("G10").Value = 154.69999999999999

# Row 11
$ws.Range("A11").Value = 46014
$ws.Range("D11").Value = 153.27000000000001
$ws.Range("E11").Value = 152.87
$ws.Range("F11").Value = 162.87
$ws.Range("G11").Value = 152.99

# Row 12
$ws.Range("A12").Value = 46014
$ws.Range("D12").Value = 153.27000000000001
$ws.Range("E12").Value = 152.87
$ws.Range("F12").Value = 162.87
$ws.Range("G12").Value = 152.99

# Row 13
$ws.Range("A13").Value = 46014
$ws.Range("D13").Value = 155.80000000000001
$ws.Range("E13").Value = 155.79
$ws.Range("F13").Value = 165.79
$ws.Range("G13").Value = 156.31

# Row 17
$ws.Range("A17").Value = 46015
$ws.Range("D17").Value = 159.21
$ws.Range("E17").Value = 160.13999999999999
$ws.Range("F17").Value = 170.14

# Row 18
$ws.Range("A18").Value = 46014
$ws.Range("D18").Value = 159.69
$ws.Range("E18").Value = 160.99
$ws.Range("F18").Value = 171

# Row 22
$ws.Range("A22").Value = 46015
$ws.Range("D22").Value = 153.93
$ws.Range("E22").Value = 153.43
$ws.Range("F22").Value = 163.03
$ws.Range("G22").Value = 154.59

# Row 23
$ws.Range("A23").Value = 46015
$ws.Range("D23").Value = 159.53
$ws.Range("E23").Value = 159.05000000000001
$ws.Range("F23").Value = 169.05

# Row 24
$ws.Range("A24").Value = 46015
$ws.Range("D24").Value = 159.31
$ws.Range("E24").Value = 159.56
$ws.Range("F24").Value = 169.56

# Row 25
$ws.Range("A25").Value = 46015
$ws.Range("D25").Value = 159.81
$ws.Range("E25").Value = 159.30000000000001
$ws.Range("F25").Value = 169.3
$ws.Range("G25").Value = 159.07

# Row 26
$ws.Range("A26").Value = 46015
$ws.Range("D26").Value = 158.74
$ws.Range("E26").Value = 160.58000000000001
$ws.Range("F26").Value = 170.58

# Row 27
$ws.Range("A27").Value = 46014
$ws.Range("D27").Value = 154.63
$ws.Range("E27").Value = 154.47999999999999
$ws.Range("F27").Value = 164.08
$ws.Range("G27").Value = 155.63999999999999

# Row 28
$ws.Range("A28").Value = 46014
$ws.Range("D28").Value = 160.27000000000001
$ws.Range("E28").Value = 159.88999999999999
$ws.Range("F28").Value = 169.89

# Row 29
$ws.Range("A29").Value = 46014
$ws.Range("D29").Value = 160.04
$ws.Range("E29").Value = 160.41
$ws.Range("F29").Value = 170.41

# Row 30
$ws.Range("A30").Value = 46014
$ws.Range("D30").Value = 160.54
$ws.Range("E30").Value = 160.16
$ws.Range("F30").Value = 170.16
$ws.Range("G30").Value = 159.93

# Row 31
$ws.Range("A31").Value = 46014
$ws.Range("D31").Value = 159.47
$ws.Range("E31").Value = 161.44999999999999
$ws.Range("F31").Value = 171.45

# Row 35
$ws.Range("A35").Value = 46015
$ws.Range("D35").Value = 152.58000000000001
$ws.Range("E35").Value = 152.63
$ws.Range("F35").Value = 161.63

# Row 36
$ws.Range("A36").Value = 46014
$ws.Range("D36").Value = 153.63999999999999
$ws.Range("E36").Value = 153.47
$ws.Range("F36").Value = 162.47

# Row 40
$ws.Range("A40").Value = 46015
$ws.Range("D40").Value = 160.37
$ws.Range("E40").Value = 160.63
$ws.Range("F40").Value = 170.63

# Row 41
$ws.Range("A41").Value = 46015
$ws.Range("D41").Value = 160.09
$ws.Range("E41").Value = 161.05000000000001
$ws.Range("F41").Value = 171.05

# Row 42
$ws.Range("A42").Value = 46014
$ws.Range("D42").Value = 160.9
$ws.Range("E42").Value = 161.54
$ws.Range("F42").Value = 171.54

# Row 43
$ws.Range("A43").Value = 46014
$ws.Range("D43").Value = 160.62
$ws.Range("E43").Value = 161.96
$ws.Range("F43").Value = 171.96

# Row 47
$ws.Range("A47").Value = 46015
$ws.Range("D47").Value = 154.16999999999999
$ws.Range("E47").Value = 155.25
$ws.Range("F47").Value = 165.25

# Row 48
$ws.Range("A48").Value = 46015
$ws.Range("D48").Value = 153.99
$ws.Range("E48").Value = 155.35
$ws.Range("F48").Value = 165.35

# Row 49
$ws.Range("A49").Value = 46014
$ws.Range("D49").Value = 154.91
$ws.Range("E49").Value = 156.03
$ws.Range("F49").Value = 166.03

# Row 50
$ws.Range("A50").Value = 46014
$ws.Range("D50").Value = 154.72999999999999
$ws.Range("E50").Value = 156.13
$ws.Range("F50").Value = 166.13

# Row 54
$ws.Range("A54").Value = 46015
$ws.Range("D54").Value = 169.98
$ws.Range("E54").Value = 169.99
$ws.Range("F54").Value = 179.99

# Row 55
$ws.Range("A55").Value = 46015
$ws.Range("D55").Value = 158.13
$ws.Range("E55").Value = 164.62
$ws.Range("F55").Value = 174.62

# Row 56
$ws.Range("A56").Value = 46015
$ws.Range("D56").Value = 159.83000000000001

# Row 57
$ws.Range("A57").Value = 46015
$ws.Range("D57").Value = 158.97999999999999
$ws.Range("E57").Value = 158.9

# Row 58
$ws.Range("A58").Value = 46015
$ws.Range("D58").Value = 154.88
$ws.Range("E58").Value = 154.94
$ws.Range("F58").Value = 164.94

# Row 59
$ws.Range("A59").Value = 46015
$ws.Range("D59").Value = 162.13999999999999
$ws.Range("E59").Value = 167.41

# Row 60
$ws.Range("A60").Value = 46014
$ws.Range("D60").Value = 170.48
$ws.Range("E60").Value = 170.95
$ws.Range("F60").Value = 180.95

# Row 61
$ws.Range("A61").Value = 46014
$ws.Range("D61").Value = 158.63999999999999
$ws.Range("E61").Value = 165.1
$ws.Range("F61").Value = 175.1

# Row 62
$ws.Range("A62").Value = 46014
$ws.Range("D62").Value = 160.01

# Row 63
$ws.Range("A63").Value = 46014
$ws.Range("D63").Value = 159.12
$ws.Range("E63").Value = 159.38

# Row 64
$ws.Range("A64").Value = 46014
$ws.Range("D64").Value = 155.02000000000001
$ws.Range("E64").Value = 155.41999999999999
$ws.Range("F64").Value = 165.42

# Row 65
$ws.Range("A65").Value = 46014
$ws.Range("D65").Value = 162.62
$ws.Range("E65").Value = 168.32
